$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume table (rows 2-51, columns B-E)
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '42.240.67'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '2.282.08'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '322.38'
$ws.Range("E5").Value = '  +2.03%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '101.79'
$ws.Range("E6").Value = '  -2.35%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  -0.80%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  -0.44%  '

$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D10").Value = '39.68'
$ws.Range("E10").Value = '  +0.02%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '8.30'
$ws.Range("E12").Value = '  -1.96%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.962'
$ws.Range("E14").Value = '  -1.18%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '15.07'
$ws.Range("E15").Value = '  -2.35%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.630.97'
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.279.18'
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '42.227.35'
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '7.34'
$ws.Range("E19").Value = '  -4.62%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0000105'
$ws.Range("E20").Value = '  -0.48%  '

$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = '12.76'
$ws.Range("E21").Value = '  +27.88%  '

$ws.Range("B22").Value = 'PancakeSwap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D22").Value = '3.65'
$ws.Range("E22").Value = '  +3.71%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '73.00'
$ws.Range("E23").Value = '  -1.05%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '268.57'
$ws.Range("E24").Value = '  -6.62%  '

$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  -3.40%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '10.83'
$ws.Range("E27").Value = '  -1.02%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '2.32'
$ws.Range("E28").Value = '  +2.14%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '22.44'
$ws.Range("E29").Value = '  -4.68%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '37.92'
$ws.Range("E30").Value = '  +5.91%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '164.04'
$ws.Range("E31").Value = '  -0.92%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '6.08'
$ws.Range("E32").Value = '  +3.17%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0876'
$ws.Range("E33").Value = '  -0.99%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '0.133'
$ws.Range("E34").Value = '  +1.09%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '2.54'
$ws.Range("E35").Value = '  -12.99%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").Value = '  -3.87%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.59'
$ws.Range("E37").Value = '  -1.33%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0354'
$ws.Range("E38").Value = '  +0.37%  '

$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '3.68'
$ws.Range("E39").Value = '  +1.57%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = '2.74'
$ws.Range("E40").Value = '  -6.96%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '1.52'
$ws.Range("E41").Value = '  +2.06%  '

$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '68.74'
$ws.Range("E42").Value = '  -2.70%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.225'
$ws.Range("E44").Value = '  -0.63%  '

$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").Value = '91.03'
$ws.Range("E45").Value = '  -10.76%  '

$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = '12.19'
$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '112.69'
$ws.Range("E47").Value = '  -3.51%  '

$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").Value = '79.82'
$ws.Range("E48").Value = '  +2.12%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '8.93'
$ws.Range("E49").Value = '  -1.93%  '

$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '5.22'
$ws.Range("E50").Value = '  -2.17%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.591.51'
$ws.Range("E51").Value = '  +1.97%  '
